$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Expected_ui_elements" column (old column C) is removed from its
# position and the columns to its right (Population_name, Expected_File_names,
# Files_to_upload, error_msg_col1, error_msg_col2, Expected_lot_options) shift
# one column to the left (D->C, E->D, F->E, G->F, H->G, I->H). Cut+Insert (as
# opposed to Delete) keeps the per-column width metadata intact as it moves.
$ws.Columns("C").Cut()
$ws.Columns("J").Insert()

# Re-add the "Expected_ui_elements" column content at the end (new column I),
# aligned one row higher than where it used to sit.
$ws.Range("I1").Value = "Expected_ui_elements"
$ws.Range("I2").Value = "Manage Line of Therapy"
$ws.Range("I3").Value = "You can view all, create new and edit or delete existing Line of Therapy from here"

# Update the error-message wording in the (now shifted) error_msg_col1 column.
$ws.Range("G2").Value = "Population filter 2 'Automation_1' is not supported"
$ws.Range("G4").Value = "Population filter 2 'Automation_2' is not supported"

# Reflect the updated view/selection state.
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("K12").Select()
